$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep Text type (matches original inlineStr cells) by
# forcing Text number format before assigning string values.
$targetRefs = @("D2","D4","D5","D6","D7","D8","E8","D9","D10","D11","D12","D13","B14","C14","D14","E14","B15","C15","D15","E15","B16","C16","D16","E16","B17","C17","D17","E17","B18","C18","D18","E18","B19","C19","D19","E19","B20","C20","D20","E20","D21","B22","C22","D22","E22","B23","C23","D23","E23","B24","C24","D24","E24","B25","C25","D25","E25","B26","C26","D26","E26","D27","D40","D41","E41","D43","D45","D47","E48")
foreach ($ref in $targetRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "242.53"
$ws.Range("D4").Value = "5.457"
$ws.Range("D5").Value = "0.05691"
$ws.Range("D6").Value = "3.415"
$ws.Range("D7").Value = "6.273"
$ws.Range("D8").Value = "1.089"
$ws.Range("E8").Value = "7FTXTokenFTTBestin24h"
$ws.Range("D9").Value = "0.8054"
$ws.Range("D10").Value = "0.1420"
$ws.Range("D11").Value = "0.07275"
$ws.Range("D12").Value = "0.03086"
$ws.Range("D13").Value = "0.03089"
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").Value = "0.1299"
$ws.Range("E14").Value = "13ProBitTokenPROB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09364"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.919"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001570"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04793"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "0.0005818"
$ws.Range("E19").Value = "18OneONE"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "0.006227"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("D21").Value = "0.0009987"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D22").Value = "0.004066"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "0.0001499"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "3.732"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "2.154"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "0.3257"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("D27").Value = "0.0003998"
$ws.Range("D40").Value = "0.03812"
$ws.Range("D41").Value = "0.006672"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D43").Value = "0.002639"
$ws.Range("D45").Value = "0.00005620"
$ws.Range("D47").Value = "0.3899"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
